# DeudoresPrueba.xlsx update:
#  - Row 27 (NEVADA, consecutivo 26): Fecha 2025-12-15 -> 2025-12-19, Valor 148700 -> 191000
#  - Row 29 (PARAÍSO MOSQUERA, consecutivo 28): Fecha 2025-12-15 -> 2025-12-19, Valor 394000 -> 362700
#  - View: selection moves to D30 (and scrolled so row 13 is at the top)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$newFecha = Get-Date -Year 2025 -Month 12 -Day 19 -Hour 0 -Minute 0 -Second 0

$ws.Range("C27").Value = $newFecha
$ws.Range("D27").Value = 191000

$ws.Range("C29").Value = $newFecha
$ws.Range("D29").Value = 362700

# Match the saved view: active cell D30, scrolled down so row 13 is at top.
$ws.Range("D30").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
